$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial date numbers, same as existing rows)
$newData = @(
    @(44308, 4, 65, 197.0712185065034),
    @(44309, 18, 69, 209.1986781069035),
    @(44310, 6, 67, 203.1349483067035),
    @(44311, 9, 58, 175.848164205803),
    @(44312, 7, 59, 178.880029105903)
)

$startRow = 234

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newData[$i]

    # Copy formatting (style) from the last existing data row (233) so the new
    # row matches the existing formatting (column A date style, borders, etc.)
    $ws.Range("A233").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $rowData[0]
    $ws.Range("B$r").Value = $rowData[1]
    $ws.Range("C$r").Value = $rowData[2]
    $ws.Range("D$r").Value = $rowData[3]
}

$excel.CutCopyMode = 0
